$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 63, pushing existing rows 63..154 down to 64..155
$ws.Rows.Item(63).Insert()

# Populate the newly inserted row 63 with the new record
$ws.Range("A63").Value = 2
$ws.Range("B63").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C63").Value = 'Coquimbo'
$ws.Range("D63").Value = 44671
$ws.Range("E63").Value = 4
$ws.Range("F63").Value = 'Fruta'
$ws.Range("G63").Value = 100109
$ws.Range("H63").Value = 'Uva'
$ws.Range("I63").Value = 100109001
$ws.Range("J63").Value = 'Uva'
$ws.Range("K63").Value = 'Red Globe'
$ws.Range("L63").Value = 'Primera'
$ws.Range("M63").Value = 400
$ws.Range("N63").Value = 7000
$ws.Range("O63").Value = 8000
$ws.Range("P63").Value = 7500
$ws.Range("Q63").Value = '$/bandeja 18 kilos'
$ws.Range("R63").Value = 'Provincia de Limarí'
$ws.Range("S63").Value = 417
$ws.Range("T63").Value = 18
